$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 22:46"

# Update Cataluña (row 5) statistics
$ws.Range("B5").Value = 7864
$ws.Range("C5").Value = 755
$ws.Range("D5").Value = 6770
$ws.Range("E5").Value = 339
